$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix spelling error: "Agents Avaialble" -> "Agents Available"
$ws.Range("G8").Value = "Agents Available"

# Apply the built-in Currency cell style/number format to the new
# column I (I2:I9), which is what the authored workbook shows
# (numFmtId 44, linked to the "Currency" named cell style).
$rng = $ws.Range("I2:I9")
$rng.Style = "Currency"

# Update the sheet's selection to reflect the newly formatted range.
$null = $ws.Range("I2:I9").Select()
